$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.249.22'
$ws.Range("E2").Value = '  +1.87%  '

$ws.Range("D3").Value = '1.796.28'
$ws.Range("E3").Value = '  +3.07%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.45'
$ws.Range("E5").Value = '  +2.01%  '

$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4570'
$ws.Range("E7").Value = '  +18.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3743'
$ws.Range("E8").Value = '  +11.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.15'
$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("E10").Value = '  +6.44%  '

$ws.Range("E11").Value = '  +4.30%  '

$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.31'
$ws.Range("E13").Value = '  +2.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.284'
$ws.Range("E14").Value = '  +3.59%  '

$ws.Range("E15").Value = '  +8.19%  '

$ws.Range("D16").Value = '1.798.81'
$ws.Range("E16").Value = '  +3.31%  '

$ws.Range("E17").Value = '  +4.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06738'
$ws.Range("E18").Value = '  +2.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.16'
$ws.Range("E19").Value = '  +3.52%  '

$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.42'
$ws.Range("E21").Value = '  +4.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.370'
$ws.Range("E22").Value = '  +3.56%  '

$ws.Range("D23").Value = '28.255.81'
$ws.Range("E23").Value = '  +1.88%  '

$ws.Range("E24").Value = '  +3.62%  '

$ws.Range("E25").Value = '  +1.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.66'
$ws.Range("E26").Value = '  +5.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.95'
$ws.Range("E27").Value = '  -1.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.350'
$ws.Range("E28").Value = '  +3.61%  '

$ws.Range("D29").Value = '2.003.06'
$ws.Range("E29").Value = '  +3.30%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.59'
$ws.Range("E30").Value = '  +4.33%  '

$ws.Range("E31").Value = '  -2.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.030'
$ws.Range("E32").Value = '  -0.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09491'
$ws.Range("E33").Value = '  +9.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.801'
$ws.Range("E34").Value = '  +1.21%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2354'
$ws.Range("E35").Value = '  +13.33%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06325'
$ws.Range("E36").Value = '  +5.01%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02342'
$ws.Range("E37").Value = '  +4.67%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.250'
$ws.Range("E38").Value = '  +3.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.02'
$ws.Range("E39").Value = '  +0.92%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6590'
$ws.Range("E40").Value = '  +3.48%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.351'
$ws.Range("E41").Value = '  +6.41%  '

$ws.Range("B42").Value = 'WEMIXTOKEN'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.485'
$ws.Range("E42").Value = '  -0.25%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.223'
$ws.Range("E43").Value = '  +3.06%  '

$ws.Range("E44").Value = '  +5.06%  '

$ws.Range("E45").Value = '  -0.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.854'
$ws.Range("E46").Value = '  +1.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6090'
$ws.Range("E47").Value = '  +3.48%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.45'
$ws.Range("E48").Value = '  +4.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.029'
$ws.Range("E49").Value = '  +3.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07144'
$ws.Range("E50").Value = '  +3.30%  '

$ws.Range("E51").Value = '  +2.14%  '
